$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.063.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.50%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.831.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.53%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.9988"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'241.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.47%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.6553"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.65%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.9997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.06%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'44.54"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +5.94%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.2938"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.39%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.07342"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.87%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'23.00"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.73%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.07668"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.63%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'1.836.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.64%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'4.983"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.39%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.6672"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.31%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'81.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.87%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'6.111"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.13%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.000008708"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +4.87%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'29.068.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.34%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'2.087.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.10%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'12.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.42%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'224.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.84%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.9999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.05%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'7.123"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.73%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.9991"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.15%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'157.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.76%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'8.499"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.13%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'0.1379"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.60%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  -0.59%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'1.506"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.17%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'4.109"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.53%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'4.016"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.19%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.201"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.06%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.05355"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.82%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("B35").Value = "'ImmutableX"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'0.7435"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.40%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("B36").Value = "'LidoDAOToken"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'1.840"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.57%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'1.157"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.03%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'2.644"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.24%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'1.298.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.15%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.01790"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.61%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'2.749"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.85%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'6.350"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +6.89%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.8987"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.45%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.9995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.29%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'103.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.04%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'1.985.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.38%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.07804"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.96%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'64.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.74%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.5138"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.46%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  -3.05%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'1.740"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.88%  "
$ws.Range("E51").Style = "Normal"
